$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '25.888.33'
    'E2' = '  -2.30%  '
    'D3' = '1.753.58'
    'E3' = '  -4.64%  '
    'D4' = '1.001'
    'E4' = '  +0.04%  '
    'D5' = '239.30'
    'E5' = '  -8.20%  '
    'D6' = '1.000'
    'E6' = '  -0.02%  '
    'D7' = '0.5091'
    'E7' = '  -5.37%  '
    'D8' = '42.32'
    'E8' = '  -5.64%  '
    'D9' = '0.2763'
    'E9' = '  -5.69%  '
    'D10' = '0.06208'
    'E10' = '  -10.39%  '
    'D11' = '1.746.98'
    'E11' = '  -5.06%  '
    'B12' = 'TRON'
    'C12' = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
    'D12' = '0.06968'
    'E12' = '  -3.04%  '
    'B13' = 'Solana'
    'C13' = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
    'D13' = '15.73'
    'E13' = '  -8.76%  '
    'D14' = '0.6134'
    'E14' = '  -15.50%  '
    'D15' = '4.533'
    'E15' = '  -8.98%  '
    'D16' = '77.53'
    'E16' = '  -12.97%  '
    'D17' = '1.001'
    'E17' = '  -0.06%  '
    'D18' = '1.000'
    'E18' = '  -0.01%  '
    'D19' = '25.901.20'
    'E19' = '  -2.33%  '
    'D20' = '0.000006944'
    'E20' = '  -11.81%  '
    'E21' = '  -15.06%  '
    'D22' = '1.968.78'
    'E22' = '  -5.48%  '
    'D23' = '4.086'
    'D24' = '5.276'
    'E24' = '  -11.86%  '
    'D25' = '8.233'
    'E25' = '  -10.37%  '
    'D26' = '138.02'
    'E26' = '  -2.66%  '
    'D27' = '1.493'
    'E27' = '  -12.40%  '
    'B28' = 'EthereumClassic'
    'C28' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D28' = '15.08'
    'E28' = '  -10.88%  '
    'B29' = 'LidoDAOToken'
    'C29' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'D29' = '1.820'
    'E29' = '  -15.89%  '
    'D30' = '103.67'
    'E30' = '  -6.57%  '
    'D31' = '0.08205'
    'E31' = '  -7.67%  '
    'D32' = '3.710'
    'E32' = '  -12.28%  '
    'D33' = '3.497'
    'E33' = '  -13.05%  '
    'D34' = '0.04541'
    'E34' = '  -6.14%  '
    'D35' = '0.9998'
    'E35' = '  +0.00%  '
    'D36' = '2.640'
    'E36' = '  -9.43%  '
    'D37' = '0.9936'
    'E37' = '  -11.94%  '
    'D38' = '0.6127'
    'E38' = '  -15.15%  '
    'D39' = '2.711'
    'E39' = '  -12.38%  '
    'D40' = '0.01558'
    'E40' = '  -8.94%  '
    'D41' = '104.13'
    'E41' = '  -2.64%  '
    'D42' = '1.001'
    'E42' = '  +0.04%  '
    'D43' = '1.896'
    'E43' = '  -17.35%  '
    'D44' = '0.3886'
    'E44' = '  -16.85%  '
    'D45' = '0.7422'
    'E45' = '  -17.74%  '
    'D46' = '4.938'
    'E46' = '  -15.69%  '
    'D47' = '0.05425'
    'E47' = '  -5.85%  '
    'D48' = '0.1118'
    'E48' = '  -10.21%  '
    'D49' = '6.015'
    'D50' = '30.14'
    'E50' = '  -13.24%  '
    'D51' = '52.84'
}

foreach ($cellref in $updates.Keys) {
    $c = $ws.Range($cellref)
    $c.NumberFormat = "@"
    $c.Value2 = $updates[$cellref]
    $c.Style = "Normal"
}
